$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5620
$ws.Range("C2").Value = 10970
$ws.Range("D2").Value = 21310
$ws.Range("E2").Value = 42120
$ws.Range("F2").Value = 83760
$ws.Range("G2").Value = 164110
$ws.Range("H2").Value = 339750
$ws.Range("I2").Value = 691640
$ws.Range("J2").Value = 1415520
